# Auto-generated edit script: swap/rotate betting-odds rows that were
# mis-ordered in the source feed for "Czech Republic 2 Liga".
# Row pairs 11/12, 19/22, 27/28, 43/44, 47/48, 108/109, 213/215 each swap
# their B:AD contents (column A "id" stays put), and rows 235/236/237
# rotate (235<-236, 236<-237, 237<-235).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Cells.Item(11, 2).Value = 6804155  # B11
$ws.Cells.Item(11, 5).Value = "FC Silon Taborsko"  # E11
$ws.Cells.Item(11, 6).Value = "Sparta Prague B"  # F11
$ws.Cells.Item(11, 9).Value = 1  # I11
$ws.Cells.Item(11, 12).Value = 2.15  # L11
$ws.Cells.Item(11, 13).Value = 3.4  # M11
$ws.Cells.Item(11, 14).Value = 2.9  # N11
$ws.Cells.Item(11, 15).Value = 1.85  # O11
$ws.Cells.Item(11, 16).Value = 3.5  # P11
$ws.Cells.Item(11, 19).Value = 1.875  # S11
$ws.Cells.Item(11, 20).Value = 1.925  # T11
$ws.Cells.Item(11, 24).Value = 0.8500000000000001  # X11
$ws.Cells.Item(11, 27).Value = 0.875  # AA11

# Row 12
$ws.Cells.Item(12, 2).Value = 6804151  # B12
$ws.Cells.Item(12, 5).Value = "Opava"  # E12
$ws.Cells.Item(12, 6).Value = "Varnsdorf"  # F12
$ws.Cells.Item(12, 9).Value = 0  # I12
$ws.Cells.Item(12, 12).Value = 1.95  # L12
$ws.Cells.Item(12, 13).Value = 3.3  # M12
$ws.Cells.Item(12, 14).Value = 3.5  # N12
$ws.Cells.Item(12, 15).Value = 1.909  # O12
$ws.Cells.Item(12, 16).Value = 3.3  # P12
$ws.Cells.Item(12, 19).Value = 1.975  # S12
$ws.Cells.Item(12, 20).Value = 1.825  # T12
$ws.Cells.Item(12, 24).Value = 0.909  # X12
$ws.Cells.Item(12, 27).Value = 0.9750000000000001  # AA12

# Row 19
$ws.Cells.Item(19, 2).Value = 6803241  # B19
$ws.Cells.Item(19, 5).Value = "Viktoria Zizkov"  # E19
$ws.Cells.Item(19, 6).Value = "SK Lisen"  # F19
$ws.Cells.Item(19, 8).Value = 1  # H19
$ws.Cells.Item(19, 9).Value = 0  # I19
$ws.Cells.Item(19, 11).Value = "H"  # K19
$ws.Cells.Item(19, 12).Value = 3  # L19
$ws.Cells.Item(19, 13).Value = 3.1  # M19
$ws.Cells.Item(19, 14).Value = 2.3  # N19
$ws.Cells.Item(19, 15).Value = 2.15  # O19
$ws.Cells.Item(19, 16).Value = 3.2  # P19
$ws.Cells.Item(19, 17).Value = 3.3  # Q19
$ws.Cells.Item(19, 19).Value = 1.85  # S19
$ws.Cells.Item(19, 20).Value = 1.95  # T19
$ws.Cells.Item(19, 21).Value = 2.25  # U19
$ws.Cells.Item(19, 22).Value = 1.775  # V19
$ws.Cells.Item(19, 23).Value = 2.025  # W19
$ws.Cells.Item(19, 24).Value = 1.15  # X19
$ws.Cells.Item(19, 25).Value = -1  # Y19
$ws.Cells.Item(19, 27).Value = 0.8500000000000001  # AA19
$ws.Cells.Item(19, 28).Value = -1  # AB19
$ws.Cells.Item(19, 29).Value = 0.7749999999999999  # AC19

# Row 22
$ws.Cells.Item(22, 2).Value = 6804159  # B22
$ws.Cells.Item(22, 5).Value = "Varnsdorf"  # E22
$ws.Cells.Item(22, 6).Value = "FC Vlasim"  # F22
$ws.Cells.Item(22, 8).Value = 3  # H22
$ws.Cells.Item(22, 9).Value = 2  # I22
$ws.Cells.Item(22, 11).Value = "D"  # K22
$ws.Cells.Item(22, 12).Value = 2.1  # L22
$ws.Cells.Item(22, 13).Value = 3.5  # M22
$ws.Cells.Item(22, 14).Value = 3  # N22
$ws.Cells.Item(22, 15).Value = 2.25  # O22
$ws.Cells.Item(22, 16).Value = 3.5  # P22
$ws.Cells.Item(22, 17).Value = 2.7  # Q22
$ws.Cells.Item(22, 19).Value = 2.025  # S22
$ws.Cells.Item(22, 20).Value = 1.775  # T22
$ws.Cells.Item(22, 21).Value = 3  # U22
$ws.Cells.Item(22, 22).Value = 1.85  # V22
$ws.Cells.Item(22, 23).Value = 1.95  # W22
$ws.Cells.Item(22, 24).Value = -1  # X22
$ws.Cells.Item(22, 25).Value = 2.5  # Y22
$ws.Cells.Item(22, 27).Value = -0.5  # AA22
$ws.Cells.Item(22, 28).Value = 0.3875  # AB22
$ws.Cells.Item(22, 29).Value = 0.8500000000000001  # AC22

# Row 27
$ws.Cells.Item(27, 2).Value = 6803160  # B27
$ws.Cells.Item(27, 5).Value = "FC Silon Taborsko"  # E27
$ws.Cells.Item(27, 6).Value = "FC Brno"  # F27
$ws.Cells.Item(27, 8).Value = 1  # H27
$ws.Cells.Item(27, 12).Value = 2.8  # L27
$ws.Cells.Item(27, 13).Value = 3.25  # M27
$ws.Cells.Item(27, 14).Value = 2.2  # N27
$ws.Cells.Item(27, 15).Value = 2.75  # O27
$ws.Cells.Item(27, 16).Value = 3.25  # P27
$ws.Cells.Item(27, 17).Value = 2.25  # Q27
$ws.Cells.Item(27, 18).Value = 0.25  # R27
$ws.Cells.Item(27, 19).Value = 1.775  # S27
$ws.Cells.Item(27, 20).Value = 2.025  # T27
$ws.Cells.Item(27, 22).Value = 1.925  # V27
$ws.Cells.Item(27, 23).Value = 1.875  # W27
$ws.Cells.Item(27, 24).Value = 1.75  # X27
$ws.Cells.Item(27, 27).Value = 0.7749999999999999  # AA27
$ws.Cells.Item(27, 29).Value = 0.925  # AC27
$ws.Cells.Item(27, 30).Value = -1  # AD27

# Row 28
$ws.Cells.Item(28, 2).Value = 6804161  # B28
$ws.Cells.Item(28, 5).Value = "SK Lisen"  # E28
$ws.Cells.Item(28, 6).Value = "FK Pribram"  # F28
$ws.Cells.Item(28, 8).Value = 0  # H28
$ws.Cells.Item(28, 12).Value = 2.5  # L28
$ws.Cells.Item(28, 13).Value = 3  # M28
$ws.Cells.Item(28, 14).Value = 2.6  # N28
$ws.Cells.Item(28, 15).Value = 2.15  # O28
$ws.Cells.Item(28, 16).Value = 3.1  # P28
$ws.Cells.Item(28, 17).Value = 3.1  # Q28
$ws.Cells.Item(28, 18).Value = -0.25  # R28
$ws.Cells.Item(28, 19).Value = 1.925  # S28
$ws.Cells.Item(28, 20).Value = 1.875  # T28
$ws.Cells.Item(28, 22).Value = 1.95  # V28
$ws.Cells.Item(28, 23).Value = 1.85  # W28
$ws.Cells.Item(28, 24).Value = 1.15  # X28
$ws.Cells.Item(28, 27).Value = 0.925  # AA28
$ws.Cells.Item(28, 29).Value = -1  # AC28
$ws.Cells.Item(28, 30).Value = 0.8500000000000001  # AD28

# Row 43
$ws.Cells.Item(43, 2).Value = 6803244  # B43
$ws.Cells.Item(43, 5).Value = "Sigma Olomouc B"  # E43
$ws.Cells.Item(43, 6).Value = "Viktoria Zizkov"  # F43
$ws.Cells.Item(43, 7).Value = 0  # G43
$ws.Cells.Item(43, 9).Value = 0  # I43
$ws.Cells.Item(43, 11).Value = "A"  # K43
$ws.Cells.Item(43, 12).Value = 2.2  # L43
$ws.Cells.Item(43, 14).Value = 2.875  # N43
$ws.Cells.Item(43, 15).Value = 2.375  # O43
$ws.Cells.Item(43, 17).Value = 2.6  # Q43
$ws.Cells.Item(43, 18).Value = 0  # R43
$ws.Cells.Item(43, 19).Value = 1.8  # S43
$ws.Cells.Item(43, 20).Value = 2  # T43
$ws.Cells.Item(43, 21).Value = 2.75  # U43
$ws.Cells.Item(43, 24).Value = -1  # X43
$ws.Cells.Item(43, 26).Value = 1.6  # Z43
$ws.Cells.Item(43, 27).Value = -1  # AA43
$ws.Cells.Item(43, 28).Value = 1  # AB43
$ws.Cells.Item(43, 29).Value = -1  # AC43
$ws.Cells.Item(43, 30).Value = 0.825  # AD43

# Row 44
$ws.Cells.Item(44, 2).Value = 6804175  # B44
$ws.Cells.Item(44, 5).Value = "MFK Chrudim"  # E44
$ws.Cells.Item(44, 6).Value = "SK Prostejov"  # F44
$ws.Cells.Item(44, 7).Value = 6  # G44
$ws.Cells.Item(44, 9).Value = 2  # I44
$ws.Cells.Item(44, 11).Value = "H"  # K44
$ws.Cells.Item(44, 12).Value = 1.95  # L44
$ws.Cells.Item(44, 14).Value = 3.3  # N44
$ws.Cells.Item(44, 15).Value = 2  # O44
$ws.Cells.Item(44, 17).Value = 3.2  # Q44
$ws.Cells.Item(44, 18).Value = -0.25  # R44
$ws.Cells.Item(44, 19).Value = 1.925  # S44
$ws.Cells.Item(44, 20).Value = 1.875  # T44
$ws.Cells.Item(44, 21).Value = 2.5  # U44
$ws.Cells.Item(44, 24).Value = 1  # X44
$ws.Cells.Item(44, 26).Value = -1  # Z44
$ws.Cells.Item(44, 27).Value = 0.925  # AA44
$ws.Cells.Item(44, 28).Value = -1  # AB44
$ws.Cells.Item(44, 29).Value = 0.9750000000000001  # AC44
$ws.Cells.Item(44, 30).Value = -1  # AD44

# Row 47
$ws.Cells.Item(47, 2).Value = 6804173  # B47
$ws.Cells.Item(47, 5).Value = "FC Vlasim"  # E47
$ws.Cells.Item(47, 6).Value = "FK Pribram"  # F47
$ws.Cells.Item(47, 7).Value = 3  # G47
$ws.Cells.Item(47, 8).Value = 1  # H47
$ws.Cells.Item(47, 10).Value = 1  # J47
$ws.Cells.Item(47, 12).Value = 2.25  # L47
$ws.Cells.Item(47, 13).Value = 3.25  # M47
$ws.Cells.Item(47, 14).Value = 2.75  # N47
$ws.Cells.Item(47, 15).Value = 1.833  # O47
$ws.Cells.Item(47, 16).Value = 3.8  # P47
$ws.Cells.Item(47, 17).Value = 3.3  # Q47
$ws.Cells.Item(47, 18).Value = -0.5  # R47
$ws.Cells.Item(47, 19).Value = 1.875  # S47
$ws.Cells.Item(47, 20).Value = 1.925  # T47
$ws.Cells.Item(47, 22).Value = 1.875  # V47
$ws.Cells.Item(47, 23).Value = 1.925  # W47
$ws.Cells.Item(47, 24).Value = 0.833  # X47
$ws.Cells.Item(47, 27).Value = 0.875  # AA47
$ws.Cells.Item(47, 28).Value = -1  # AB47
$ws.Cells.Item(47, 29).Value = 0.875  # AC47
$ws.Cells.Item(47, 30).Value = -1  # AD47

# Row 48
$ws.Cells.Item(48, 2).Value = 6804172  # B48
$ws.Cells.Item(48, 5).Value = "Dukla Praha"  # E48
$ws.Cells.Item(48, 6).Value = "Vysocina Jihlava"  # F48
$ws.Cells.Item(48, 7).Value = 1  # G48
$ws.Cells.Item(48, 8).Value = 0  # H48
$ws.Cells.Item(48, 10).Value = 0  # J48
$ws.Cells.Item(48, 12).Value = 1.8  # L48
$ws.Cells.Item(48, 13).Value = 3.4  # M48
$ws.Cells.Item(48, 14).Value = 3.8  # N48
$ws.Cells.Item(48, 15).Value = 1.5  # O48
$ws.Cells.Item(48, 16).Value = 4.2  # P48
$ws.Cells.Item(48, 17).Value = 5  # Q48
$ws.Cells.Item(48, 18).Value = -1  # R48
$ws.Cells.Item(48, 19).Value = 1.825  # S48
$ws.Cells.Item(48, 20).Value = 1.975  # T48
$ws.Cells.Item(48, 22).Value = 1.925  # V48
$ws.Cells.Item(48, 23).Value = 1.875  # W48
$ws.Cells.Item(48, 24).Value = 0.5  # X48
$ws.Cells.Item(48, 27).Value = 0  # AA48
$ws.Cells.Item(48, 28).Value = 0  # AB48
$ws.Cells.Item(48, 29).Value = -1  # AC48
$ws.Cells.Item(48, 30).Value = 0.875  # AD48

# Row 108
$ws.Cells.Item(108, 2).Value = 6804221  # B108
$ws.Cells.Item(108, 5).Value = "MFK Vyskov"  # E108
$ws.Cells.Item(108, 6).Value = "Opava"  # F108
$ws.Cells.Item(108, 7).Value = 1  # G108
$ws.Cells.Item(108, 8).Value = 3  # H108
$ws.Cells.Item(108, 9).Value = 1  # I108
$ws.Cells.Item(108, 10).Value = 2  # J108
$ws.Cells.Item(108, 11).Value = "A"  # K108
$ws.Cells.Item(108, 12).Value = 1.95  # L108
$ws.Cells.Item(108, 19).Value = 1.8  # S108
$ws.Cells.Item(108, 20).Value = 2  # T108
$ws.Cells.Item(108, 22).Value = 1.95  # V108
$ws.Cells.Item(108, 23).Value = 1.85  # W108
$ws.Cells.Item(108, 24).Value = -1  # X108
$ws.Cells.Item(108, 26).Value = 3.2  # Z108
$ws.Cells.Item(108, 27).Value = -1  # AA108
$ws.Cells.Item(108, 28).Value = 1  # AB108
$ws.Cells.Item(108, 29).Value = 0.95  # AC108

# Row 109
$ws.Cells.Item(109, 2).Value = 6803338  # B109
$ws.Cells.Item(109, 5).Value = "MFK Chrudim"  # E109
$ws.Cells.Item(109, 6).Value = "Hanacka Slavia Kromeriz"  # F109
$ws.Cells.Item(109, 7).Value = 3  # G109
$ws.Cells.Item(109, 8).Value = 1  # H109
$ws.Cells.Item(109, 9).Value = 3  # I109
$ws.Cells.Item(109, 10).Value = 0  # J109
$ws.Cells.Item(109, 11).Value = "H"  # K109
$ws.Cells.Item(109, 12).Value = 1.909  # L109
$ws.Cells.Item(109, 19).Value = 1.825  # S109
$ws.Cells.Item(109, 20).Value = 1.975  # T109
$ws.Cells.Item(109, 22).Value = 1.825  # V109
$ws.Cells.Item(109, 23).Value = 1.975  # W109
$ws.Cells.Item(109, 24).Value = 0.75  # X109
$ws.Cells.Item(109, 26).Value = -1  # Z109
$ws.Cells.Item(109, 27).Value = 0.825  # AA109
$ws.Cells.Item(109, 28).Value = -1  # AB109
$ws.Cells.Item(109, 29).Value = 0.825  # AC109

# Row 213
$ws.Cells.Item(213, 2).Value = 6804289  # B213
$ws.Cells.Item(213, 5).Value = "MFK Chrudim"  # E213
$ws.Cells.Item(213, 6).Value = "Dukla Praha"  # F213
$ws.Cells.Item(213, 7).Value = 0  # G213
$ws.Cells.Item(213, 8).Value = 2  # H213
$ws.Cells.Item(213, 9).Value = 0  # I213
$ws.Cells.Item(213, 11).Value = "A"  # K213
$ws.Cells.Item(213, 12).Value = 3.1  # L213
$ws.Cells.Item(213, 13).Value = 3.4  # M213
$ws.Cells.Item(213, 14).Value = 2  # N213
$ws.Cells.Item(213, 15).Value = 3.6  # O213
$ws.Cells.Item(213, 16).Value = 3.6  # P213
$ws.Cells.Item(213, 17).Value = 1.8  # Q213
$ws.Cells.Item(213, 18).Value = 0.5  # R213
$ws.Cells.Item(213, 19).Value = 1.975  # S213
$ws.Cells.Item(213, 20).Value = 1.825  # T213
$ws.Cells.Item(213, 21).Value = 2.75  # U213
$ws.Cells.Item(213, 22).Value = 1.95  # V213
$ws.Cells.Item(213, 23).Value = 1.85  # W213
$ws.Cells.Item(213, 24).Value = -1  # X213
$ws.Cells.Item(213, 26).Value = 0.8  # Z213
$ws.Cells.Item(213, 27).Value = -1  # AA213
$ws.Cells.Item(213, 28).Value = 0.825  # AB213
$ws.Cells.Item(213, 29).Value = -1  # AC213
$ws.Cells.Item(213, 30).Value = 0.8500000000000001  # AD213

# Row 215
$ws.Cells.Item(215, 2).Value = 6804286  # B215
$ws.Cells.Item(215, 5).Value = "FC Silon Taborsko"  # E215
$ws.Cells.Item(215, 6).Value = "Opava"  # F215
$ws.Cells.Item(215, 7).Value = 2  # G215
$ws.Cells.Item(215, 8).Value = 1  # H215
$ws.Cells.Item(215, 9).Value = 1  # I215
$ws.Cells.Item(215, 11).Value = "H"  # K215
$ws.Cells.Item(215, 12).Value = 1.727  # L215
$ws.Cells.Item(215, 13).Value = 3.5  # M215
$ws.Cells.Item(215, 14).Value = 4  # N215
$ws.Cells.Item(215, 15).Value = 1.8  # O215
$ws.Cells.Item(215, 16).Value = 3.3  # P215
$ws.Cells.Item(215, 17).Value = 3.8  # Q215
$ws.Cells.Item(215, 18).Value = -0.5  # R215
$ws.Cells.Item(215, 19).Value = 1.875  # S215
$ws.Cells.Item(215, 20).Value = 1.925  # T215
$ws.Cells.Item(215, 21).Value = 2.25  # U215
$ws.Cells.Item(215, 22).Value = 1.875  # V215
$ws.Cells.Item(215, 23).Value = 1.925  # W215
$ws.Cells.Item(215, 24).Value = 0.8  # X215
$ws.Cells.Item(215, 26).Value = -1  # Z215
$ws.Cells.Item(215, 27).Value = 0.875  # AA215
$ws.Cells.Item(215, 28).Value = -1  # AB215
$ws.Cells.Item(215, 29).Value = 0.875  # AC215
$ws.Cells.Item(215, 30).Value = -1  # AD215

# Row 235
$ws.Cells.Item(235, 2).Value = 6836420  # B235
$ws.Cells.Item(235, 5).Value = "FK Pribram"  # E235
$ws.Cells.Item(235, 6).Value = "MFK Vyskov"  # F235
$ws.Cells.Item(235, 7).Value = 2  # G235
$ws.Cells.Item(235, 9).Value = 1  # I235
$ws.Cells.Item(235, 12).Value = 2  # L235
$ws.Cells.Item(235, 13).Value = 3.4  # M235
$ws.Cells.Item(235, 14).Value = 3.1  # N235
$ws.Cells.Item(235, 15).Value = 1.833  # O235
$ws.Cells.Item(235, 16).Value = 4  # P235
$ws.Cells.Item(235, 17).Value = 3.25  # Q235
$ws.Cells.Item(235, 18).Value = -0.5  # R235
$ws.Cells.Item(235, 19).Value = 1.85  # S235
$ws.Cells.Item(235, 20).Value = 1.95  # T235
$ws.Cells.Item(235, 21).Value = 2.75  # U235
$ws.Cells.Item(235, 22).Value = 1.8  # V235
$ws.Cells.Item(235, 23).Value = 2  # W235
$ws.Cells.Item(235, 24).Value = 0.833  # X235
$ws.Cells.Item(235, 27).Value = 0.8500000000000001  # AA235
$ws.Cells.Item(235, 29).Value = -1  # AC235
$ws.Cells.Item(235, 30).Value = 1  # AD235

# Row 236
$ws.Cells.Item(236, 2).Value = 6920870  # B236
$ws.Cells.Item(236, 5).Value = "SK Lisen"  # E236
$ws.Cells.Item(236, 6).Value = "MFK Chrudim"  # F236
$ws.Cells.Item(236, 7).Value = 1  # G236
$ws.Cells.Item(236, 8).Value = 2  # H236
$ws.Cells.Item(236, 9).Value = 0  # I236
$ws.Cells.Item(236, 11).Value = "A"  # K236
$ws.Cells.Item(236, 12).Value = 2.1  # L236
$ws.Cells.Item(236, 13).Value = 3.5  # M236
$ws.Cells.Item(236, 14).Value = 2.875  # N236
$ws.Cells.Item(236, 15).Value = 2.1  # O236
$ws.Cells.Item(236, 16).Value = 3.5  # P236
$ws.Cells.Item(236, 17).Value = 2.875  # Q236
$ws.Cells.Item(236, 18).Value = -0.25  # R236
$ws.Cells.Item(236, 19).Value = 1.9  # S236
$ws.Cells.Item(236, 20).Value = 1.9  # T236
$ws.Cells.Item(236, 21).Value = 2.5  # U236
$ws.Cells.Item(236, 22).Value = 1.95  # V236
$ws.Cells.Item(236, 23).Value = 1.85  # W236
$ws.Cells.Item(236, 24).Value = -1  # X236
$ws.Cells.Item(236, 26).Value = 1.875  # Z236
$ws.Cells.Item(236, 27).Value = -1  # AA236
$ws.Cells.Item(236, 28).Value = 0.8999999999999999  # AB236
$ws.Cells.Item(236, 29).Value = 0.95  # AC236
$ws.Cells.Item(236, 30).Value = -1  # AD236

# Row 237
$ws.Cells.Item(237, 2).Value = 6920869  # B237
$ws.Cells.Item(237, 5).Value = "Dukla Praha"  # E237
$ws.Cells.Item(237, 6).Value = "Sparta Prague B"  # F237
$ws.Cells.Item(237, 7).Value = 6  # G237
$ws.Cells.Item(237, 8).Value = 0  # H237
$ws.Cells.Item(237, 9).Value = 2  # I237
$ws.Cells.Item(237, 11).Value = "H"  # K237
$ws.Cells.Item(237, 12).Value = 1.571  # L237
$ws.Cells.Item(237, 13).Value = 4  # M237
$ws.Cells.Item(237, 14).Value = 4.333  # N237
$ws.Cells.Item(237, 15).Value = 1.4  # O237
$ws.Cells.Item(237, 16).Value = 4.75  # P237
$ws.Cells.Item(237, 17).Value = 5.5  # Q237
$ws.Cells.Item(237, 18).Value = -1.25  # R237
$ws.Cells.Item(237, 19).Value = 1.875  # S237
$ws.Cells.Item(237, 20).Value = 1.925  # T237
$ws.Cells.Item(237, 21).Value = 3.25  # U237
$ws.Cells.Item(237, 22).Value = 1.825  # V237
$ws.Cells.Item(237, 23).Value = 1.975  # W237
$ws.Cells.Item(237, 24).Value = 0.3999999999999999  # X237
$ws.Cells.Item(237, 26).Value = -1  # Z237
$ws.Cells.Item(237, 27).Value = 0.875  # AA237
$ws.Cells.Item(237, 28).Value = -1  # AB237
$ws.Cells.Item(237, 29).Value = 0.825  # AC237

